# Replace each two-digit/one-digit division answer with its updated value.
# Every original text run is unique within the document, so a plain literal
# Find/Replace (MatchWildcards=$false) targeting each exact old string is safe.
$d = $word.ActiveDocument

$d.Content.Find.Execute("54÷6=9, 0", $true, $false, $false, $false, $false, $true, 1, $false, "96÷4=24, 0", 2) | Out-Null
$d.Content.Find.Execute("47÷7=6, 5", $true, $false, $false, $false, $false, $true, 1, $false, "49÷9=5, 4", 2) | Out-Null
$d.Content.Find.Execute("94÷3=31, 1", $true, $false, $false, $false, $false, $true, 1, $false, "19÷6=3, 1", 2) | Out-Null
$d.Content.Find.Execute("94÷2=47, 0", $true, $false, $false, $false, $false, $true, 1, $false, "98÷5=19, 3", 2) | Out-Null
$d.Content.Find.Execute("82÷6=13, 4", $true, $false, $false, $false, $false, $true, 1, $false, "88÷9=9, 7", 2) | Out-Null
$d.Content.Find.Execute("39÷5=7, 4", $true, $false, $false, $false, $false, $true, 1, $false, "43÷8=5, 3", 2) | Out-Null
$d.Content.Find.Execute("56÷5=11, 1", $true, $false, $false, $false, $false, $true, 1, $false, "99÷5=19, 4", 2) | Out-Null
$d.Content.Find.Execute("54÷9=6, 0", $true, $false, $false, $false, $false, $true, 1, $false, "15÷6=2, 3", 2) | Out-Null
$d.Content.Find.Execute("48÷5=9, 3", $true, $false, $false, $false, $false, $true, 1, $false, "23÷9=2, 5", 2) | Out-Null
$d.Content.Find.Execute("74÷8=9, 2", $true, $false, $false, $false, $false, $true, 1, $false, "56÷6=9, 2", 2) | Out-Null
$d.Content.Find.Execute("45÷8=5, 5", $true, $false, $false, $false, $false, $true, 1, $false, "55÷4=13, 3", 2) | Out-Null
$d.Content.Find.Execute("68÷5=13, 3", $true, $false, $false, $false, $false, $true, 1, $false, "80÷8=10, 0", 2) | Out-Null
$d.Content.Find.Execute("59÷2=29, 1", $true, $false, $false, $false, $false, $true, 1, $false, "79÷8=9, 7", 2) | Out-Null
$d.Content.Find.Execute("43÷2=21, 1", $true, $false, $false, $false, $false, $true, 1, $false, "28÷5=5, 3", 2) | Out-Null
$d.Content.Find.Execute("58÷6=9, 4", $true, $false, $false, $false, $false, $true, 1, $false, "48÷4=12, 0", 2) | Out-Null
$d.Content.Find.Execute("88÷2=44, 0", $true, $false, $false, $false, $false, $true, 1, $false, "25÷9=2, 7", 2) | Out-Null
$d.Content.Find.Execute("62÷4=15, 2", $true, $false, $false, $false, $false, $true, 1, $false, "32÷8=4, 0", 2) | Out-Null
$d.Content.Find.Execute("75÷6=12, 3", $true, $false, $false, $false, $false, $true, 1, $false, "61÷5=12, 1", 2) | Out-Null
$d.Content.Find.Execute("22÷4=5, 2", $true, $false, $false, $false, $false, $true, 1, $false, "30÷5=6, 0", 2) | Out-Null
$d.Content.Find.Execute("52÷2=26, 0", $true, $false, $false, $false, $false, $true, 1, $false, "92÷9=10, 2", 2) | Out-Null
$d.Content.Find.Execute("11÷5=2, 1", $true, $false, $false, $false, $false, $true, 1, $false, "55÷6=9, 1", 2) | Out-Null
$d.Content.Find.Execute("26÷3=8, 2", $true, $false, $false, $false, $false, $true, 1, $false, "86÷3=28, 2", 2) | Out-Null
$d.Content.Find.Execute("53÷3=17, 2", $true, $false, $false, $false, $false, $true, 1, $false, "40÷6=6, 4", 2) | Out-Null
$d.Content.Find.Execute("73÷7=10, 3", $true, $false, $false, $false, $false, $true, 1, $false, "44÷7=6, 2", 2) | Out-Null
$d.Content.Find.Execute("81÷8=10, 1", $true, $false, $false, $false, $false, $true, 1, $false, "39÷4=9, 3", 2) | Out-Null
